$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend row 5 with the two missing trailing columns (X5, Y5) ---
$ws.Range("X5").Value = 0.6499990000000011
$ws.Range("Y5").Value = "Up"

# --- Append a brand-new row 6 (next day's scan results) ---
$ws.Range("A6").Value = 42647.884305555555
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = "Random"
$ws.Range("Q6").Value = 46.357611069683557
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.0591
$ws.Range("S6").NumberFormat = "0.00%"
$ws.Range("T6").Value = -0.0421
$ws.Range("T6").NumberFormat = "0.00%"
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = 0
